# Add VTA lines to transit capacity files
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transitLineToVehicle")

# Insert 3 new rows right before the existing "137_A" row (currently row
# 492), pushing the HSR rows + trailing header-echo row down to 495-498.
# Excel carries each shifted row's own formatting along with it, so row
# 495 (the old row 492, "137_A") keeps its "first row of a vehicle group"
# top-border style (6 / 11) that we want to reuse for the new group below.
$ws.Rows("492:494").Insert()

# Row 492: 111_LRCCW
$ws.Range("A492").Value = "111_LRCCW"
$ws.Range("B492").Value = "VTA LRT"
$ws.Range("D492").Value = "CCW"
$ws.Range("E492").Value = "GREEN LINE"
$ws.Range("F492").Value = "LRV2"
$ws.Range("G492").Value = "LRV2"
$ws.Range("H492").Value = "LRV2"
$ws.Range("C492").Formula = '=RIGHT($A492,LEN($A492)-FIND("_",$A492))'

# Row 493: 111_902LRT
$ws.Range("A493").Value = "111_902LRT"
$ws.Range("B493").Value = "VTA LRT"
$ws.Range("D493").Value = 902
$ws.Range("E493").Value = "GREEN LINE"
$ws.Range("F493").Value = "LRV2"
$ws.Range("G493").Value = "LRV2"
$ws.Range("H493").Value = "LRV2"
$ws.Range("C493").Formula = '=RIGHT($A493,LEN($A493)-FIND("_",$A493))'

# Row 494: 111_LRTWCC
$ws.Range("A494").Value = "111_LRTWCC"
$ws.Range("B494").Value = "VTA LRT"
$ws.Range("D494").Value = "WCC"
$ws.Range("E494").Value = "GREEN LINE"
$ws.Range("F494").Value = "LRV2"
$ws.Range("G494").Value = "LRV2"
$ws.Range("H494").Value = "LRV2"
$ws.Range("C494").Formula = '=RIGHT($A494,LEN($A494)-FIND("_",$A494))'

# Row 492 is the first row of the new VTA group, so (like every other
# first-row-of-group in this table, e.g. row 495 "137_A") it carries the
# "group divider" top-border style. Row 495 already has exactly that
# style (it's the shifted-down copy of the original row 492), so copy
# formats from there onto the new row 492.
$ws.Range("A495:H495").Copy() | Out-Null
$ws.Range("A492:H492").PasteSpecial(-4122) | Out-Null

# Column B keeps the same top-border style for all three new VTA rows
# (matches the source workbook, where the "VTA LRT" label was filled
# down with its original formatting intact).
$ws.Range("B492").Copy() | Out-Null
$ws.Range("B493").PasteSpecial(-4122) | Out-Null
$ws.Range("B494").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# View state: scroll position / selection moved along with the new rows.
$ws.Application.ActiveWindow.ScrollRow = 471
$ws.Range("E494").Select()

$ws3 = $wb.Worksheets.Item("transitVehicleToCapacity")
$ws3.Application.ActiveWindow.ScrollRow = 14
